$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range stays text-formatted (the source sheet stores all
# values - including "82%"-style percentages and "2.4"-style decimals - as
# literal text, not numbers) while we overwrite the updated cells below.
$dataRange = $ws.Range("A2:G21")
$dataRange.NumberFormat = "@"

$changes = @(
    ,@("C2", "1.3")
    ,@("D2", "6.2")
    ,@("E2", "85%")
    ,@("F2", "58%")
    ,@("G2", "2.94")
    ,@("C3", "2.3")
    ,@("D3", "4.7")
    ,@("E3", "72%")
    ,@("F3", "63%")
    ,@("G3", "3.12")
    ,@("E4", "60%")
    ,@("F4", "42%")
    ,@("G4", "2.21")
    ,@("C5", "2.3")
    ,@("D5", "4.2")
    ,@("E5", "73%")
    ,@("F5", "46%")
    ,@("G5", "2.24")
    ,@("C6", "2.8")
    ,@("D6", "4.3")
    ,@("E6", "78%")
    ,@("F6", "59%")
    ,@("G6", "2.97")
    ,@("C7", "2.2")
    ,@("E7", "87%")
    ,@("F7", "66%")
    ,@("G7", "3.00")
    ,@("B8", "Lazio")
    ,@("C8", "2.5")
    ,@("D8", "5.4")
    ,@("E8", "66%")
    ,@("F8", "42%")
    ,@("G8", "2.33")
    ,@("B9", "Napoli")
    ,@("C9", "2.1")
    ,@("D9", "6.4")
    ,@("E9", "79%")
    ,@("F9", "58%")
    ,@("G9", "2.76")
    ,@("B10", "Fiorentina")
    ,@("C10", "2.2")
    ,@("D10", "5.2")
    ,@("E10", "69%")
    ,@("F10", "44%")
    ,@("G10", "2.53")
    ,@("B11", "Torino")
    ,@("C11", "2.1")
    ,@("D11", "4.8")
    ,@("E11", "61%")
    ,@("F11", "27%")
    ,@("G11", "1.82")
    ,@("B12", "Monza")
    ,@("C12", "2.4")
    ,@("F12", "39%")
    ,@("G12", "2.36")
    ,@("C13", "2.2")
    ,@("D13", "3.9")
    ,@("E13", "67%")
    ,@("F13", "39%")
    ,@("G13", "2.27")
    ,@("B14", "Lecce")
    ,@("C14", "2.4")
    ,@("D14", "4.6")
    ,@("E14", "73%")
    ,@("F14", "42%")
    ,@("G14", "2.36")
    ,@("B15", "Cagliari")
    ,@("C15", "2.0")
    ,@("D15", "4.7")
    ,@("E15", "82%")
    ,@("F15", "57%")
    ,@("G15", "2.79")
    ,@("B16", "Hellas Verona")
    ,@("C16", "2.3")
    ,@("D16", "3.6")
    ,@("E16", "61%")
    ,@("F16", "46%")
    ,@("G16", "2.27")
    ,@("B17", "Empoli")
    ,@("C17", "2.2")
    ,@("D17", "4.7")
    ,@("E17", "61%")
    ,@("F17", "33%")
    ,@("G17", "2.24")
    ,@("B18", "Udinese")
    ,@("C18", "2.5")
    ,@("D18", "4.2")
    ,@("E18", "78%")
    ,@("F18", "41%")
    ,@("G18", "2.44")
    ,@("C19", "1.8")
    ,@("E19", "85%")
    ,@("F19", "64%")
    ,@("G19", "3.12")
    ,@("C20", "1.8")
    ,@("D20", "5.5")
    ,@("F20", "58%")
    ,@("G20", "3.15")
    ,@("C21", "2.4")
    ,@("D21", "4.0")
    ,@("E21", "85%")
    ,@("G21", "2.91")
)

foreach ($change in $changes) {
    $ws.Range($change[0]).Value = $change[1]
}

# Restore the default (General/Normal) style so the cells match the
# original workbook's look (plain shared-string cells, no explicit style).
$dataRange.Style = "Normal"

Write-Output "Applied $($changes.Count) cell updates"
